$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45566.82485457524
$ws.Cells.Item(2, 2).Value = 1000
$ws.Cells.Item(2, 3).Value = 48.51067352294922
$ws.Cells.Item(2, 4).Value = 20.61402011924086
$ws.Cells.Item(2, 5).Value = 1000
$ws.Cells.Item(2, 6).Value = 1983.274856799464
$ws.Cells.Item(2, 7).Value = 983.274856799464

$ws.Cells.Item(3, 1).Value = 45597.82485457524
$ws.Cells.Item(3, 2).Value = 1000
$ws.Cells.Item(3, 3).Value = 49.56356048583984
$ws.Cells.Item(3, 4).Value = 20.17611305962768
$ws.Cells.Item(3, 5).Value = 2000
$ws.Cells.Item(3, 6).Value = 3924.41867579446
$ws.Cells.Item(3, 7).Value = 1924.41867579446

$ws.Cells.Item(4, 1).Value = 45627.82485457524
$ws.Cells.Item(4, 2).Value = 1000
$ws.Cells.Item(4, 3).Value = 46.4828987121582
$ws.Cells.Item(4, 4).Value = 21.51328827817782
$ws.Cells.Item(4, 5).Value = 3000
$ws.Cells.Item(4, 6).Value = 5994.212121341945
$ws.Cells.Item(4, 7).Value = 2994.212121341945

$ws.Cells.Item(5, 1).Value = 45658.82485457524
$ws.Cells.Item(5, 2).Value = 1000
$ws.Cells.Item(5, 3).Value = 42.75
$ws.Cells.Item(5, 4).Value = 23.39181286549708
$ws.Cells.Item(5, 5).Value = 4000
$ws.Cells.Item(5, 6).Value = 8244.738415715574
$ws.Cells.Item(5, 7).Value = 4244.738415715574

$ws.Cells.Item(6, 1).Value = 45689.82485457524
$ws.Cells.Item(6, 2).Value = 6338
$ws.Cells.Item(6, 3).Value = 48.36999893188477
$ws.Cells.Item(6, 4).Value = 131.0316340698136
$ws.Cells.Item(6, 5).Value = 10338
$ws.Cells.Item(6, 6).Value = 20851.2918096093
$ws.Cells.Item(6, 7).Value = 10513.2918096093

$ws.Cells.Item(7, 1).Value = 45717.82485457524
$ws.Cells.Item(7, 2).Value = 1000
$ws.Cells.Item(7, 3).Value = 48.65999984741211
$ws.Cells.Item(7, 4).Value = 20.55076044257701
$ws.Cells.Item(7, 5).Value = 11338
$ws.Cells.Item(7, 6).Value = 22828.48045297485
$ws.Cells.Item(7, 7).Value = 11490.48045297485

$ws.Cells.Item(8, 1).Value = 45748.82485457524
$ws.Cells.Item(8, 2).Value = 1000
$ws.Cells.Item(8, 3).Value = 56.90999984741211
$ws.Cells.Item(8, 4).Value = 17.57160433458468
$ws.Cells.Item(8, 5).Value = 12338
$ws.Cells.Item(8, 6).Value = 24519.04448991796
$ws.Cells.Item(8, 7).Value = 12181.04448991796

$ws.Cells.Item(9, 1).Value = 45778.82485457524
$ws.Cells.Item(9, 2).Value = 1000
$ws.Cells.Item(9, 3).Value = 58.72999954223633
$ws.Cells.Item(9, 4).Value = 17.02707317885877
$ws.Cells.Item(9, 5).Value = 13338
$ws.Cells.Item(9, 6).Value = 26157.21918486721
$ws.Cells.Item(9, 7).Value = 12819.21918486721

$ws.Cells.Item(10, 1).Value = 45809.82485457524
$ws.Cells.Item(10, 2).Value = 1000
$ws.Cells.Item(10, 3).Value = 65.31999969482422
$ws.Cells.Item(10, 4).Value = 15.30924685658315
$ws.Cells.Item(10, 5).Value = 14338
$ws.Cells.Item(10, 6).Value = 27630.12181092304
$ws.Cells.Item(10, 7).Value = 13292.12181092304

$ws.Cells.Item(11, 1).Value = 45839.82485457524
$ws.Cells.Item(11, 2).Value = 1779
$ws.Cells.Item(11, 3).Value = 67.58000183105469
$ws.Cells.Item(11, 4).Value = 26.32435560518889
$ws.Cells.Item(11, 5).Value = 16117
$ws.Cells.Item(11, 6).Value = 30162.7880395976
$ws.Cells.Item(11, 7).Value = 14045.7880395976

$ws.Cells.Item(12, 1).Value = 45870.82485457524
$ws.Cells.Item(12, 2).Value = 1000
$ws.Cells.Item(12, 3).Value = 65.02999877929688
$ws.Cells.Item(12, 4).Value = 15.37751835724104
$ws.Cells.Item(12, 5).Value = 17117
$ws.Cells.Item(12, 6).Value = 31642.25906666922
$ws.Cells.Item(12, 7).Value = 14525.25906666922

$ws.Cells.Item(13, 1).Value = 45901.82485457524
$ws.Cells.Item(13, 2).Value = 1000
$ws.Cells.Item(13, 3).Value = 80.05000305175781
$ws.Cells.Item(13, 4).Value = 12.492191903521
$ws.Cells.Item(13, 5).Value = 18117
$ws.Cells.Item(13, 6).Value = 32844.13283827003
$ws.Cells.Item(13, 7).Value = 14727.13283827003
